# Applies the "Reset the Values to null" change-request edit to the CR
# document: merges a handful of split runs back together (no visible text
# change), fixes a couple of real text edits (Yangmin -> Huang Cun, Yang Min
# -> Huang Cun), relocates the "_GoBack" bookmark, and retargets the two
# embedded-picture runs' language tag.

$d = $word.ActiveDocument

# --- 1. Remove the stray _GoBack bookmark that sat right after
#        "CHANGE DESCRIPTION:" (it gets re-added further down, see step 5).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. "Yangmin" (flagged by the spell checker) -> " Huang Cun"
$d.Content.Find.Execute("Yangmin", $false, $false, $false, $false, $false, `
                         $true, 1, $false, " Huang Cun", 2)

# --- 3. "Yang Min     " (5 trailing spaces) -> "Huang Cun  " (2 trailing spaces)
$d.Content.Find.Execute("Yang Min     ", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Huang Cun  ", 2)

# --- 4. The "__________" placeholder on the Coding/Testing man-days line
#        becomes "___1______" (a bookmark "_GoBack" now sits between the
#        "1" and the trailing underscores).
$d.Content.Find.Execute("2 man-days" + [char]9 + "__________", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "2 man-days" + [char]9 + "___1______", 2)

# Re-insert the _GoBack bookmark between "___1" and "______" on that line.
$r = $d.Content
$r.Find.Execute("___1______", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r.Find.Found) {
    $gb = $d.Range($r.Start + 4, $r.Start + 4)
    $d.Bookmarks.Add("_GoBack", $gb)
}

$word.Selection.HomeKey(6)
